# Refresh the cryptocurrency price / volume snapshot (GitHub Actions run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "26.701.83"
$ws.Range("E2").Value = "  +0.18%  "
# Row 3
$ws.Range("D3").Value = "1.599.51"
$ws.Range("E3").Value = "  +0.01%  "
# Row 4
$ws.Range("E4").Value = "  +0.17%  "
# Row 5
$ws.Range("D5").Value = "'211.34"
$ws.Range("E5").Value = "  -0.07%  "
# Row 6
$ws.Range("E6").Value = "  -0.92%  "
# Row 7
$ws.Range("E7").Value = "  +0.15%  "
# Row 8
$ws.Range("E8").Value = "  +0.34%  "
# Row 9
$ws.Range("E9").Value = "  +0.95%  "
# Row 10
$ws.Range("E10").Value = "  +0.14%  "
# Row 11
$ws.Range("E11").Value = "  +0.78%  "
# Row 12
$ws.Range("D12").Value = "1.823.69"
$ws.Range("E12").Value = "  -0.02%  "
# Row 13
$ws.Range("D13").Value = "1.597.54"
$ws.Range("E13").Value = "  -1.43%  "
# Row 15
$ws.Range("E15").Value = "  +0.34%  "
# Row 16
$ws.Range("D16").Value = "'65.39"
$ws.Range("E16").Value = "  +0.97%  "
# Row 17
$ws.Range("D17").Value = "26.674.90"
$ws.Range("E17").Value = "  +0.12%  "
# Row 18
$ws.Range("E18").Value = "  +3.66%  "
# Row 19
$ws.Range("D19").Value = "'210.20"
$ws.Range("E19").Value = "  +0.94%  "
# Row 20
$ws.Range("E20").Value = "  +0.23%  "
# Row 21
$ws.Range("D21").Value = "'7.17"
$ws.Range("E21").Value = "  +3.23%  "
# Row 22
$ws.Range("E22").Value = "  +0.39%  "
# Row 23
$ws.Range("E23").Value = "  -0.07%  "
# Row 24
$ws.Range("D24").Value = "'8.93"
$ws.Range("E24").Value = "  +0.95%  "
# Row 25
$ws.Range("D25").Value = "'143.17"
$ws.Range("E25").Value = "  -1.61%  "
# Row 26
$ws.Range("E26").Value = "  +0.10%  "
# Row 27
$ws.Range("E27").Value = "  -0.20%  "
# Row 28
$ws.Range("E28").Value = "  +0.04%  "
# Row 29
$ws.Range("E29").Value = "  +0.26%  "
# Row 31
$ws.Range("E31").Value = "  -0.04%  "
# Row 32
$ws.Range("E32").Value = "  +0.53%  "
# Row 33
$ws.Range("E33").Value = "  +1.65%  "
# Row 34
$ws.Range("D34").Value = "1.290.42"
$ws.Range("E34").Value = "  +0.68%  "
# Row 35
$ws.Range("E35").Value = "  -5.54%  "
# Row 36
$ws.Range("E36").Value = "  +0.98%  "
# Row 37
$ws.Range("E37").Value = "  +0.30%  "
# Row 38
$ws.Range("E38").Value = "  -0.27%  "
# Row 39
$ws.Range("E39").Value = "  +16.60%  "
# Row 40
$ws.Range("E40").Value = "  -1.99%  "
# Row 41
$ws.Range("E41").Value = "  -0.59%  "
# Row 42
$ws.Range("E42").Value = "  -0.08%  "
# Row 43
$ws.Range("E43").Value = "  -0.86%  "
# Row 44
$ws.Range("D44").Value = "'63.22"
$ws.Range("E44").Value = "  -1.13%  "
# Row 45
$ws.Range("D45").Value = "1.727.83"
$ws.Range("E45").Value = "  -0.50%  "
# Row 46
$ws.Range("D46").Value = "'91.24"
$ws.Range("E46").Value = "  +1.65%  "
# Row 47
$ws.Range("D47").Value = "'1.57"
$ws.Range("E47").Value = "  -1.37%  "
# Row 48
$ws.Range("E48").Value = "  -1.02%  "
# Row 49
$ws.Range("E49").Value = "  +0.54%  "
# Row 50
$ws.Range("E50").Value = "  +0.08%  "
# Row 51
$ws.Range("D51").Value = "'7.35"
$ws.Range("E51").Value = "  -1.57%  "
